$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: LOGIN06 / error_user
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$ws.Rows.Item(7).RowHeight = 17.25
$ws.Range("A7").Value = "LOGIN06"
$ws.Range("B7").Value = "error_user"

# New row 8: LOGIN07 / visual_user
$ws.Range("A6:C6").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Rows.Item(8).RowHeight = 17.25
$ws.Range("A8").Value = "LOGIN07"
$ws.Range("B8").Value = "visual_user"

$excel.CutCopyMode = 0
